# Refresh the cryptocurrency table with the latest scraped price /
# volume(1h) figures, plus the rows whose coin ranking shifted
# position (TheGraph/dogwifhat/Kaspa and Stellar/CoreDAO).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.092.53"
$ws.Range("E2").Value = "'  +2.39%  "
$ws.Range("D3").Value = "'3.179.94"
$ws.Range("E3").Value = "'  +1.40%  "
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'594.55"
$ws.Range("E5").Value = "'  +4.20%  "
$ws.Range("D6").Value = "'153.26"
$ws.Range("D8").Value = "'3.177.73"
$ws.Range("E8").Value = "'  +1.28%  "
$ws.Range("E9").Value = "'  +3.60%  "
$ws.Range("E10").Value = "'  +1.08%  "
$ws.Range("D11").Value = "'6.01"
$ws.Range("E11").Value = "'  -0.29%  "
$ws.Range("D12").Value = "'0.514"
$ws.Range("E12").Value = "'  +4.44%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "'  +1.60%  "
$ws.Range("D14").Value = "'38.84"
$ws.Range("E14").Value = "'  +5.63%  "
$ws.Range("D15").Value = "'3.708.98"
$ws.Range("E15").Value = "'  +1.40%  "
$ws.Range("D16").Value = "'66.124.61"
$ws.Range("E16").Value = "'  +2.06%  "
$ws.Range("D17").Value = "'7.40"
$ws.Range("E17").Value = "'  +5.27%  "
$ws.Range("D18").Value = "'3.188.77"
$ws.Range("E18").Value = "'  +1.40%  "
$ws.Range("E19").Value = "'  +0.89%  "
$ws.Range("D20").Value = "'507.16"
$ws.Range("E20").Value = "'  +1.77%  "
$ws.Range("D21").Value = "'15.25"
$ws.Range("E21").Value = "'  +3.88%  "
$ws.Range("D22").Value = "'0.733"
$ws.Range("E22").Value = "'  +3.36%  "
$ws.Range("D23").Value = "'7.99"
$ws.Range("E23").Value = "'  +4.61%  "
$ws.Range("D24").Value = "'15.02"
$ws.Range("E24").Value = "'  -0.51%  "
$ws.Range("D25").Value = "'84.61"
$ws.Range("E25").Value = "'  +1.24%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "'  +0.09%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "'  +4.22%  "
$ws.Range("D28").Value = "'2.98"
$ws.Range("E28").Value = "'  +3.94%  "
$ws.Range("D29").Value = "'2.28"
$ws.Range("E29").Value = "'  +6.54%  "
$ws.Range("D30").Value = "'6.99"
$ws.Range("E30").Value = "'  +14.50%  "
$ws.Range("D31").Value = "'2.87"
$ws.Range("E31").Value = "'  +4.29%  "
$ws.Range("D32").Value = "'28.02"
$ws.Range("E32").Value = "'  +2.58%  "
$ws.Range("E33").Value = "'  +3.15%  "
$ws.Range("E34").Value = "'  +0.01%  "
$ws.Range("D35").Value = "'6.47"
$ws.Range("E35").Value = "'  +0.94%  "
$ws.Range("D36").Value = "'54.69"
$ws.Range("E36").Value = "'  +0.62%  "
$ws.Range("D37").Value = "'483.69"
$ws.Range("E37").Value = "'  +4.47%  "
$ws.Range("D38").Value = "'0.0894"
$ws.Range("E38").Value = "'  +0.28%  "
$ws.Range("D39").Value = "'0.0419"
$ws.Range("E39").Value = "'  +1.29%  "
$ws.Range("D40").Value = "'8.85"
$ws.Range("E40").Value = "'  +3.35%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.121"
$ws.Range("E41").Value = "'  +4.87%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.299"
$ws.Range("E42").Value = "'  +6.72%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.82"
$ws.Range("E43").Value = "'  -4.01%  "
$ws.Range("D44").Value = "'0.0₃0652"
$ws.Range("E44").Value = "'  +14.21%  "
$ws.Range("D45").Value = "'2.897.84"
$ws.Range("E45").Value = "'  -4.09%  "
$ws.Range("D46").Value = "'2.40"
$ws.Range("E46").Value = "'  +0.15%  "
$ws.Range("D47").Value = "'28.28"
$ws.Range("E47").Value = "'  +0.88%  "
$ws.Range("B49").Value = "CoreDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D49").Value = "'2.67"
$ws.Range("E49").Value = "'  +12.47%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.116"
$ws.Range("E50").Value = "'  +2.53%  "
$ws.Range("D51").Value = "'2.31"
$ws.Range("E51").Value = "'  +3.86%  "
